$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Accented characters used in client names (Í, Ñ) built from code points
# to avoid literal-encoding issues for non-ASCII bytes in the script source.
$cI = [char]0x00CD
$cN = [char]0x00D1

# Clear the previous data rows (A2:E33) before writing the refreshed data set.
$ws.Range("A2:E33").ClearContents()

# Final data: Consecutivo, Cliente, Fecha (serial), Valor, Pagado
$data = @(
    @(1, ("ALISO"), 45996, 196000, 0),
    @(2, ("ARROZ PAISA SUBA"), 46001, 166000, 0),
    @(3, ("CAMILIN"), 45997, 166000, 0),
    @(4, ("CAMPO VERDE TOCANCIPA"), 46002, 540000, 0),
    @(5, ("CAMPO VERDE ZIPAUIRA"), 46002, 681800, 0),
    @(6, ("CARNILANDIA"), 46000, 436000, 0),
    @(7, ("CARNIVOROS"), 46001, 499000, 0),
    @(8, ("CIMARRON DORADO"), 46000, 473300, 0),
    @(9, ("CLIENTE PAOLA"), 46000, 92000, 0),
    @(10, ("COCINA CHINA"), 46003, 170000, 0),
    @(11, ("COCINA CHINA"), 45998, 170000, 0),
    @(12, ("DARWIN FUTBOL"), 45921, 200000, 0),
    @(13, ("DAVIDCITO"), 45947, 100000, 0),
    @(14, ("El CEBU"), 45947, 181800, 0),
    @(15, ("EL RUBY"), 45992, 85100, 0),
    @(16, ("FRANCO"), 45996, 20000, 0),
    @(17, ("FRANCO"), 46003, 600400, 0),
    @(18, ("LA SELECTA"), 45912, 82000, 0),
    @(19, ("LOS PAISANOS"), 46002, 262500, 0),
    @(20, ("MERKA FRUVER ALEJANDRO"), 46002, 388100, 0),
    @(21, ("MERKA FRUVER DEXI"), 45995, 339000, 0),
    @(22, ("MERKA FRUVER DEXI"), 45988, 15400, 0),
    @(23, ("MULTICARNEA"), 46003, 912400, 0),
    @(24, ("NEVADA"), 45996, 229000, 0),
    @(25, ("NEVADA"), 46000, 164000, 0),
    @(26, ("NOVILLON SAN MATEO"), 45971, 83000, 0),
    @(27, ("PARA" + $cI + "SO FUNZA"), 45996, 202000, 0),
    @(28, ("PARA" + $cI + "SO MOSQUERQ"), 46003, 300000, 0),
    @(29, ("PINILLA"), 45924, 16000, 0),
    @(30, ("PINILLA"), 45931, 166000, 0),
    @(31, ("PLACITA MADRILE" + $cN + "A"), 46003, 100000, 0),
    @(32, ("PLAZA JESSICA"), 46000, 1238000, 0),
    @(33, ("PLAZA JESSICA"), 45999, 971300, 0),
    @(34, ("PORTAL ZIPA"), 46002, 664000, 0),
    @(35, ("PUNTA DE ANCA"), 46000, 7600, 0),
    @(36, ("SAMY"), 46003, 92300, 0),
    @(37, ("SANDRA 20 DE JULIO"), 46000, 300000, 0),
    @(38, ("SANTANDER SUR"), 45993, 80000, 0),
    @(39, ("SANTANDER SUR"), 45997, 250700, 0),
    @(40, ("VNZLNO PUNTA ANCA"), 45992, 82000, 0)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = [bool]($row[4])
    $r++
}

# Re-apply the date format across the full refreshed range so the
# newly added rows (34-41, beyond the original 33-row extent) pick up
# the same date style as the rest of column C instead of "General".
$ws.Range("C2:C41").NumberFormat = "YYYY-MM-DD"

# Column widths added in the refreshed layout (values tuned so the
# engine's character-width/pixel quantization lands on the closest
# achievable stored width to the authored 19.140625 / 10.42578125).
$ws.Columns.Item(2).ColumnWidth = 18.3
$ws.Columns.Item(3).ColumnWidth = 9.6

# Restore the saved selection/active cell.
$ws.Range("G36").Select()
